$d = $word.ActiveDocument

# Locate the paragraph that currently reads "2.  " and grab the first of the
# two blank paragraphs that immediately follow it (this is the insertion
# point described by the diff: a new "a)"/"b)" block goes in right after the
# "2." line, ahead of the pre-existing trailing blank paragraphs).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "2.  ") {
        $target = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($null -eq $target) {
    throw "Could not locate the '2.  ' paragraph to anchor the new content after."
}

$xml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t xml:space="preserve">a) Constraints are animals will eat each other if left alone in wrong pair </w:t></w:r></w:p>
<w:p/>
<w:p><w:r><w:t xml:space="preserve">b) </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>sub</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> goals are to find right pairing of animals and seed </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>withing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pairs in boat </w:t></w:r></w:p>
<w:p/>
<w:p/>
<w:p/>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$target.Range.InsertXML($xml) | Out-Null
